$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 107 (shifts existing rows 107:201 down to 108:202,
# expanding the used range from A1:R201 to A1:R202).
$ws.Rows.Item(107).Insert()

# Populate the newly inserted row 107 with the new weekly record.
$ws.Cells.Item(107, 1).Value = 6
$ws.Cells.Item(107, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(107, 3).Value = 'Metropolitana'
$ws.Cells.Item(107, 4).Value = 44669
$ws.Cells.Item(107, 5).Value = 13
$ws.Cells.Item(107, 6).Value = 100112001
$ws.Cells.Item(107, 7).Value = 'Berenjena'
$ws.Cells.Item(107, 8).Value = 'Sin especificar'
$ws.Cells.Item(107, 9).Value = 'Primera'
$ws.Cells.Item(107, 10).Value = 600
$ws.Cells.Item(107, 11).Value = 4000
$ws.Cells.Item(107, 12).Value = 5000
$ws.Cells.Item(107, 13).Value = 4417
$ws.Cells.Item(107, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(107, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(107, 16).Value = 88
$ws.Cells.Item(107, 17).Value = 50
$ws.Cells.Item(107, 18).Value = 'Hortaliza'
